$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")
$ws.Range("A1").Value = "TEST"
